# feat: add 2022-Q3 data
#
# 1. Insert a new "2022-Q3" worksheet (fund holdings detail) right after the
#    "总计" summary sheet and before the existing "2021-Q4" sheet.
# 2. Update the "总计" summary sheet with a new leading row for 2022-Q3 and
#    push the existing quarters down (2020-Q4 becomes the new last row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create + position the new "2022-Q3" worksheet
# ---------------------------------------------------------------------------
# NOTE: Move() invalidates previously-captured worksheet references in this
# host, so every handle used below is (re)fetched *after* the move completes.
$tmp = $wb.Worksheets.Add()
$tmp.Name = "2022-Q3"
$moveBefore = $wb.Worksheets.Item("2021-Q4")
$tmp.Move($moveBefore)

$newSheet = $wb.Worksheets.Item("2022-Q3")
$template = $wb.Worksheets.Item("2021-Q4")

# Match the look of the other quarter sheets: bold/boxed header row (B1:H1)
# and the boxed index column (A2:A4), copied from the existing "2021-Q4"
# sheet which shares the identical layout.
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2:A4").Copy()
$newSheet.Range("A2:A4").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$newSheet.Range("A4").Value = 2

# Row 2 - 003513 / 中邮消费升级灵活配置混合
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "003513"
$newSheet.Range("C2").Value = "中邮消费升级灵活配置混合"
$newSheet.Range("D2").Value = "0.56"
$newSheet.Range("E2").Value = "30.56"
$newSheet.Range("F2").Value = "1.64"
$newSheet.Range("G2").Value = "0.0092"
$newSheet.Range("H2").Value = 9

# Row 3 - 006231 / 国融融君灵活配置混合A
$newSheet.Range("B3:G3").NumberFormat = "@"
$newSheet.Range("B3").Value = "006231"
$newSheet.Range("C3").Value = "国融融君灵活配置混合A"
$newSheet.Range("D3").Value = "0.10"
$newSheet.Range("E3").Value = "55.44"
$newSheet.Range("F3").Value = "2.48"
$newSheet.Range("G3").Value = "0.0025"
$newSheet.Range("H3").Value = 2

# Row 4 - 006232 / 国融融君灵活配置混合C (note: holding value is a plain 0,
# not a formatted text string, matching the source data)
$newSheet.Range("B4:F4").NumberFormat = "@"
$newSheet.Range("B4").Value = "006232"
$newSheet.Range("C4").Value = "国融融君灵活配置混合C"
$newSheet.Range("D4").Value = "0.00"
$newSheet.Range("E4").Value = "55.44"
$newSheet.Range("F4").Value = "2.48"
$newSheet.Range("G4").Value = 0
$newSheet.Range("H4").Value = 2

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Snapshot the existing quarter rows (2021-Q4 .. 2020-Q4) before they shift.
$existingRows = @()
for ($r = 2; $r -le 6; $r++) {
  $existingRows += , @(
    $summary.Cells.Item($r, 2).Value2,
    $summary.Cells.Item($r, 3).Value2,
    $summary.Cells.Item($r, 4).Value2
  )
}

# Extend the boxed index-column style (column A) down to the new last row.
$summary.Range("A6").Copy()
$summary.Range("A7").PasteSpecial(-4122)

# New first data row: 2022-Q3 totals.
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.01

# Shift all previously-existing rows down by one.
for ($i = 0; $i -lt $existingRows.Count; $i++) {
  $targetRow = 3 + $i
  $summary.Cells.Item($targetRow, 2).Value = $existingRows[$i][0]
  $summary.Cells.Item($targetRow, 3).Value = $existingRows[$i][1]
  $summary.Cells.Item($targetRow, 4).Value = $existingRows[$i][2]
}

# Column A is just a running index (0,1,2,...) - extend it to the new row.
$summary.Range("A7").Value = 5

# Restore the originally-active sheet/selection (creating/moving sheets above
# shifts the active tab as a side effect).
$summary.Activate()
[void]$summary.Range("A1").Select()
